$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four values in row 3 from 1 to 0.99
$ws.Range("T3").Value = 0.99
$ws.Range("U3").Value = 0.99
$ws.Range("V3").Value = 0.99
$ws.Range("W3").Value = 0.99

# Update the selection to be the single cell T3 (matches active cell in diff)
$ws.Range("T3").Select()
